# "Made Timeslot code working"
# - Rebuilds new_dataset (sheet4) with the new Flight/Gate/Cost/Gate_com/arr_time/dep_time
#   layout (4 gate options per flight, 20 rows total) instead of the old
#   Flight/Gate(text)/arr_time/dep_time layout.
# - Adds a new "Sheet2" worksheet after "new_dataset" that holds the same
#   header/column layout, but only the first data row per flight is filled in.
# - Leaves a few cosmetic selection / active-cell tweaks matching the new
#   layout.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 0. Cosmetic: Model sheet selection moves to C2:C21
#    (done first so the later sheet activations below "win" and leave
#    Sheet2 as the active/visible tab, matching the saved workbook)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Model")
$ws1.Range("C2:C21").Select()

# ---------------------------------------------------------------------------
# 1. Rewrite "new_dataset" (Flight, Gate, Cost, Gate_com, arr_time, dep_time)
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("new_dataset")
$ws4.Cells.Clear()

$ws4.Cells.Item(1,1).Value = "Flight"
$ws4.Cells.Item(1,2).Value = "Gate"
$ws4.Cells.Item(1,3).Value = "Cost"
$ws4.Cells.Item(1,4).Value = "Gate_com"
$ws4.Cells.Item(1,5).Value = "arr_time"
$ws4.Cells.Item(1,6).Value = "dep_time"

$datasetRows = @(
    @(1,1,1,1,0.38541666666666669,0.41666666666666669),
    @(1,2,2,1,0.38541666666666669,0.41666666666666669),
    @(1,3,3,1,0.38541666666666669,0.41666666666666669),
    @(1,4,3,1,0.38541666666666669,0.41666666666666669),
    @(2,1,2,1,0.40625,0.46875),
    @(2,2,1,1,0.40625,0.46875),
    @(2,3,4,1,0.40625,0.46875),
    @(2,4,3,1,0.40625,0.46875),
    @(3,1,1,1,0.45833333333333331,0.5),
    @(3,2,3,1,0.45833333333333331,0.5),
    @(3,3,4,1,0.45833333333333331,0.5),
    @(3,4,2,1,0.45833333333333331,0.5),
    @(4,1,1,1,0.45833333333333331,0.51041666666666663),
    @(4,2,3,1,0.45833333333333331,0.51041666666666663),
    @(4,3,4,1,0.45833333333333331,0.51041666666666663),
    @(4,4,2,1,0.45833333333333331,0.51041666666666663),
    @(5,1,1,1,0.48958333333333331,0.54166666666666663),
    @(5,2,3,1,0.48958333333333331,0.54166666666666663),
    @(5,3,1,1,0.48958333333333331,0.54166666666666663),
    @(5,4,3,1,0.48958333333333331,0.54166666666666663)
)

$r = 2
foreach ($row in $datasetRows) {
    $ws4.Cells.Item($r,1).Value = $row[0]
    $ws4.Cells.Item($r,2).Value = $row[1]
    $ws4.Cells.Item($r,3).Value = $row[2]
    $ws4.Cells.Item($r,4).Value = $row[3]
    $ws4.Cells.Item($r,5).Value = $row[4]
    $ws4.Cells.Item($r,5).NumberFormat = "h:mm"
    $ws4.Cells.Item($r,6).Value = $row[5]
    $ws4.Cells.Item($r,6).NumberFormat = "h:mm"
    $r += 1
}

$ws4.Range("A1:F21").Select()

# ---------------------------------------------------------------------------
# 2. Add "Sheet2" right after "new_dataset" with the same header layout, but
#    only the first gate option of every flight filled in (rows 2,6,10,14,18)
# ---------------------------------------------------------------------------
$ws5 = $wb.Worksheets.Add($null, $ws4)
$ws5.Name = "Sheet2"

$ws5.Cells.Item(1,1).Value = "Flight"
$ws5.Cells.Item(1,2).Value = "Gate"
$ws5.Cells.Item(1,3).Value = "Cost"
$ws5.Cells.Item(1,4).Value = "Gate_com"
$ws5.Cells.Item(1,5).Value = "arr_time"
$ws5.Cells.Item(1,6).Value = "dep_time"

$sheet2Rows = @{
    2  = @(1,1,1,1,0.38541666666666669,0.41666666666666669)
    6  = @(2,1,2,1,0.40625,0.46875)
    10 = @(3,1,1,1,0.45833333333333331,0.5)
    14 = @(4,1,1,1,0.45833333333333331,0.51041666666666663)
    18 = @(5,1,1,1,0.48958333333333331,0.54166666666666663)
}

for ($row = 2; $row -le 21; $row++) {
    if ($sheet2Rows.ContainsKey($row)) {
        $vals = $sheet2Rows[$row]
        $ws5.Cells.Item($row,1).Value = $vals[0]
        $ws5.Cells.Item($row,2).Value = $vals[1]
        $ws5.Cells.Item($row,3).Value = $vals[2]
        $ws5.Cells.Item($row,4).Value = $vals[3]
        $ws5.Cells.Item($row,5).Value = $vals[4]
        $ws5.Cells.Item($row,6).Value = $vals[5]
    }
    $ws5.Cells.Item($row,5).NumberFormat = "h:mm"
    $ws5.Cells.Item($row,6).NumberFormat = "h:mm"
}

$ws5.Activate()
$ws5.Cells.Item(4,3).Select()
